$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update staff names in column A to include Mr./Mrs. honorific prefixes
# (rows whose names already carried a "Dr." title are left untouched).
# NOTE: assignment order below matches the original author's edit order
# so that newly-created shared-string entries land in the same sequence.
$ws.Range("A3").Value  = "Mrs. LEENA JASMINE J S"
$ws.Range("A8").Value  = "Mrs. PONNRAJAKUMARI M"
$ws.Range("A9").Value  = "Mr. SARAVANAN V"
$ws.Range("A10").Value = "Mrs. SUMATHI S"
$ws.Range("A11").Value = "Mr. MAGESH V"
$ws.Range("A5").Value  = "Mrs. THILAGAM K"
$ws.Range("A12").Value = "Mr. GANGADURAI E"
$ws.Range("A13").Value = "MR. SUBRAMANIAN G"
$ws.Range("A14").Value = "Mrs. SHANKARI R"
$ws.Range("A15").Value = "Mr. THEIVANATHAN G"
$ws.Range("A16").Value = "Mr. RADHAKRISHNAN K"
$ws.Range("A17").Value = "Mr. GNANA ARUN JOHNSON "
$ws.Range("A18").Value = "Mrs. DOLLY IRENE J"
$ws.Range("A19").Value = "Mrs. SARUPRIYA S"
$ws.Range("A20").Value = "Mrs. DEEPA N"
$ws.Range("A21").Value = "Mrs. KAVITHA S"
$ws.Range("A22").Value = "Mrs. SUJATHA R"
$ws.Range("A23").Value = "Mr. SESHAIAH M A"
$ws.Range("A24").Value = "Mrs. NANDHINI M"
$ws.Range("A25").Value = "Mrs. RENUKA P"
$ws.Range("A26").Value = "Mrs. SUGANTHI AMUDHAN DR"
$ws.Range("A27").Value = "Mrs. ARUL THILAGAVATHI M"
$ws.Range("A28").Value = "Mrs. PRAVEENA V"
$ws.Range("A29").Value = "Mrs. SUSANNA P"
$ws.Range("A30").Value = "Mrs. VINODHINI R DR"
$ws.Range("A31").Value = "Mrs. RAMYA T"
$ws.Range("A32").Value = "Mrs. MANIMEGALAI M"
$ws.Range("A33").Value = "Mrs. DHANALAKSHMI G"
$ws.Range("A34").Value = "Mr. SUBRAMANIAN N"

# Move the view/selection to match the updated worksheet position
$ws.Range("A34").Select()
